# Updates the "cryptos" price table (Sheet1) to the latest scrape:
#  - refreshes Price (D) / Volume(1h) (E) figures for every existing row
#  - inserts a new "BitDAO" row at position 20, shifting Avalanche..Cronos
#    down by one row and dropping NEARProtocol off the bottom of the list
#  - updates RenderToken's (row 51) Price/Volume in place
#
# Price values are forced to text ("@" number format) before being written
# because several of them are numeric-looking strings (e.g. "0.5010",
# "0.06350", "14.93") that Excel would otherwise silently reinterpret as
# numbers and strip the significant trailing zeros from.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.986.67"
$ws.Range("E2").Value = "  -2.33%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.815.92"
$ws.Range("E3").Value = "  -1.41%  "

$ws.Range("E4").Value = "  -1.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.56"
$ws.Range("E5").Value = "  -2.58%  "

$ws.Range("E6").Value = "  -1.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4205"
$ws.Range("E7").Value = "  -2.32%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3663"
$ws.Range("E8").Value = "  -1.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07197"
$ws.Range("E9").Value = "  -1.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8401"
$ws.Range("E10").Value = "  -3.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.73"
$ws.Range("E11").Value = "  -3.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.813.71"
$ws.Range("E12").Value = "  -1.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.620"
$ws.Range("E13").Value = "  -1.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07062"
$ws.Range("E14").Value = "  -0.73%  "

$ws.Range("E15").Value = "  -3.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.82"
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("E17").Value = "  -1.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008778"
$ws.Range("E18").Value = "  -2.06%  "

$ws.Range("E19").Value = "  -1.02%  "

$ws.Range("B20").Value = "BitDAO"
$ws.Range("C20").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.5010"
$ws.Range("E20").Value = "  -3.02%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.93"
$ws.Range("E21").Value = "  -3.40%  "

$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.031.47"
$ws.Range("E22").Value = "  -2.19%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.108"
$ws.Range("E23").Value = "  -1.87%  "

$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.79"
$ws.Range("E24").Value = "  -2.56%  "

$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.041.84"
$ws.Range("E25").Value = "  -1.55%  "

$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.973"
$ws.Range("E26").Value = "  -1.62%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.56"
$ws.Range("E27").Value = "  -2.61%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.216"
$ws.Range("E28").Value = "  +2.54%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.25"
$ws.Range("E29").Value = "  -1.77%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.201"
$ws.Range("E30").Value = "  -2.91%  "

$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "115.84"
$ws.Range("E31").Value = "  -2.35%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08772"
$ws.Range("E32").Value = "  -1.85%  "

$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.171"
$ws.Range("E33").Value = "  -4.58%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.955"
$ws.Range("E34").Value = "  +2.69%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7362"
$ws.Range("E35").Value = "  -5.03%  "

$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.401"
$ws.Range("E36").Value = "  -3.04%  "

$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  -1.15%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.087"
$ws.Range("E38").Value = "  -4.10%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01951"
$ws.Range("E39").Value = "  -0.98%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05215"
$ws.Range("E40").Value = "  -2.12%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.250"
$ws.Range("E41").Value = "  -0.52%  "

$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.864"
$ws.Range("E42").Value = "  -2.05%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1683"
$ws.Range("E43").Value = "  +0.09%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5021"
$ws.Range("E44").Value = "  -1.67%  "

$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.555"
$ws.Range("E45").Value = "  -2.62%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.54"
$ws.Range("E46").Value = "  -1.80%  "

$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.97"
$ws.Range("E47").Value = "  -2.86%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4721"
$ws.Range("E48").Value = "  -0.26%  "

$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("E49").Value = "  -1.15%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06350"
$ws.Range("E50").Value = "  -1.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.877"
$ws.Range("E51").Value = "  +2.07%  "
